$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge split runs that contain the same logical sentence/phrase back into
#    a single run by doing a literal Find/Replace over the full text. Word's
#    find-replace collapses the matched span into one run (taking the
#    formatting of the first matched run), which also removes the old
#    _GoBack bookmark markers that previously sat inside the Foursquare
#    paragraph.
# ---------------------------------------------------------------------------

function Merge-Text($text) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $null = $find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

Merge-Text("Data Requirements and Source")
Merge-Text("Foursquare API - By using this API, we will get all the venues in each neighbourhood")
Merge-Text("Geocoder python - Neighbourhood coordinates of each locality of Bangalore city")

# ---------------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark: it used to sit inside the Foursquare
#    bullet (right after "API"); now it should start at the very beginning
#    of the document and end right after the Geocoder bullet.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$geocoderPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Geocoder python*") {
        $geocoderPara = $p
    }
}

$bmRange = $d.Range(0, $geocoderPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3. Swap the Georgia font for Helvetica everywhere (ascii/hAnsi only - the
#    eastAsia/cs fallback to Times New Roman is left untouched), and give the
#    bold heading run the Helvetica font too. This has to run paragraph by
#    paragraph because Font.Name assignment on a genuinely empty paragraph
#    (just the end-of-paragraph mark) silently fails to stick in this
#    engine; working around that by temporarily inserting a character,
#    stamping the font, and removing the character again while keeping the
#    mark's formatting.
# ---------------------------------------------------------------------------

$paragraphs = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paragraphs += $d.Paragraphs.Item($i)
}

foreach ($p in $paragraphs) {
    $r = $p.Range
    if ($r.Text.Length -gt 0) {
        $r.Font.Name = "Helvetica"
    } else {
        $r.InsertBefore("x")
        $p.Range.Font.Name = "Helvetica"
        $tmp = $d.Range($p.Range.Start, $p.Range.Start + 1)
        $tmp.Delete()
    }
}
